# Refresh cached market-price-derived columns (H:N) on each job sheet.
# Source data (currentAveragePrice*, LevePrice*, LeveProfit*) is pulled from an
# external market-board snapshot and pasted in as static values - no formulas
# live on these cells, so each row below is a straight value overwrite.
$wb = $excel.ActiveWorkbook

# ALC!row 17 - One for the Road
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 501.80554
$ws.Range("J17").Value = 371.48572
$ws.Range("L17").Value = 1114.45716
$ws.Range("N17").Value = -1450.45716

# ALC!row 38 - Just Give Him a Serum
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1498.6428
$ws.Range("I38").Value = 90.083336
$ws.Range("J38").Value = 9950
$ws.Range("K38").Value = 270.250008
$ws.Range("L38").Value = 29850
$ws.Range("M38").Value = 101.749992
$ws.Range("N38").Value = -30594

# ALC!row 70 - Consecrating Congregation
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 13973539
$ws.Range("I70").Value = 33534374
$ws.Range("J70").Value = 1514.2858
$ws.Range("K70").Value = 100603122
$ws.Range("L70").Value = 4542.857400000001
$ws.Range("M70").Value = -100602852
$ws.Range("N70").Value = -5082.857400000001

# ALC!row 73 - Curbing the Contagion (L)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 13973539
$ws.Range("I73").Value = 33534374
$ws.Range("J73").Value = 1514.2858
$ws.Range("K73").Value = 100603122
$ws.Range("L73").Value = 4542.857400000001
$ws.Range("M73").Value = -100602186
$ws.Range("N73").Value = -6414.857400000001

# ALC!row 112 - Making Ends Meet
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2285.4211
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 2320.162
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 6960.485999999999
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -9176.485999999999

# ALC!row 129 - Practical Command
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1158.3833
$ws.Range("J129").Value = 1182.862
$ws.Range("L129").Value = 3548.586
$ws.Range("N129").Value = -13548.586

# ALC!row 137 - Cutting Edge of Culinary Quality
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2452.054
$ws.Range("I137").Value = 1503.9445
$ws.Range("J137").Value = 3350.2632
$ws.Range("K137").Value = 4511.833500000001
$ws.Range("L137").Value = 10050.7896
$ws.Range("M137").Value = -1961.833500000001
$ws.Range("N137").Value = -15150.7896

# ALC!row 138 - All-night Crafting
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3530.7585
$ws.Range("I138").Value = 1939.4814
$ws.Range("J138").Value = 4246.8335
$ws.Range("K138").Value = 5818.4442
$ws.Range("L138").Value = 12740.5005
$ws.Range("M138").Value = -678.4441999999999
$ws.Range("N138").Value = -23020.5005

# ALC!row 141 - Remedy for Reason
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 5608.1924
$ws.Range("I141").Value = 1832.52
$ws.Range("K141").Value = 5497.559999999999
$ws.Range("M141").Value = -317.5599999999995

# ARM!row 32 - Ingot We Trust
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15364.98
$ws.Range("I32").Value = 15062.0625
$ws.Range("K32").Value = 15062.0625
$ws.Range("M32").Value = -14775.0625

# ARM!row 61 - Dealing with the Tough Stuff
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2406.5881
$ws.Range("I61").Value = 2333.913
$ws.Range("J61").Value = 2558.5454
$ws.Range("K61").Value = 2333.913
$ws.Range("L61").Value = 2558.5454
$ws.Range("M61").Value = -2121.913
$ws.Range("N61").Value = -2982.5454

# ARM!row 122 - Haste for High Durium
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 6259.6924
$ws.Range("I122").Value = 6393.7144
$ws.Range("J122").Value = 5696.8
$ws.Range("K122").Value = 19181.1432
$ws.Range("L122").Value = 17090.4
$ws.Range("M122").Value = -16731.1432
$ws.Range("N122").Value = -21990.4

# ARM!row 132 - Don't Bore Me, Ore Me
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6060.303
$ws.Range("I132").Value = 7166.2383
$ws.Range("J132").Value = 4124.9165
$ws.Range("K132").Value = 21498.7149
$ws.Range("L132").Value = 12374.7495
$ws.Range("M132").Value = -18968.7149
$ws.Range("N132").Value = -17434.7495

# ARM!row 136 - Metal with Mettle
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2406.5881
$ws.Range("I136").Value = 2333.913
$ws.Range("J136").Value = 2558.5454
$ws.Range("K136").Value = 7001.739
$ws.Range("L136").Value = 7675.6362
$ws.Range("M136").Value = -4451.739
$ws.Range("N136").Value = -12775.6362

# BSM!row 105 - Ingot to Wing It
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 11907914
$ws.Range("I105").Value = 15875827
$ws.Range("J105").Value = 4173.6665
$ws.Range("K105").Value = 15875827
$ws.Range("L105").Value = 4173.6665
$ws.Range("M105").Value = -15874080
$ws.Range("N105").Value = -7667.6665

# CUL!row 68 - Such a Butter Face
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 103156.766
$ws.Range("I68").Value = 179365.75
$ws.Range("J68").Value = 1544.7858
$ws.Range("K68").Value = 538097.25
$ws.Range("L68").Value = 4634.357400000001
$ws.Range("M68").Value = -537286.25
$ws.Range("N68").Value = -6256.357400000001

# CUL!row 71 - No Margarine of Error (L)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 103156.766
$ws.Range("I71").Value = 179365.75
$ws.Range("J71").Value = 1544.7858
$ws.Range("K71").Value = 1614291.75
$ws.Range("L71").Value = 13903.0722
$ws.Range("M71").Value = -1610235.75
$ws.Range("N71").Value = -22015.0722

# CUL!row 122 - Salt of the North
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 706.7857
$ws.Range("I122").Value = 616.3333
$ws.Range("J122").Value = 1249.5
$ws.Range("K122").Value = 5546.9997
$ws.Range("L122").Value = 11245.5
$ws.Range("M122").Value = -3096.9997
$ws.Range("N122").Value = -16145.5

# CUL!row 133 - Friends Are Food
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 4752.154
$ws.Range("I133").Value = 1803.3334
$ws.Range("J133").Value = 7279.7144
$ws.Range("K133").Value = 5410.0002
$ws.Range("L133").Value = 21839.1432
$ws.Range("M133").Value = -350.0002000000004
$ws.Range("N133").Value = -31959.1432

# GSM!row 43 - Get the Green Stuff
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5500
$ws.Range("I43").Value = 5500
$ws.Range("J43").Value = 5500
$ws.Range("K43").Value = 5500
$ws.Range("L43").Value = 5500
$ws.Range("M43").Value = -5349
$ws.Range("N43").Value = -5802

# GSM!row 122 - Awarding Academic Excellence
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3411.6365
$ws.Range("I122").Value = 3461.2856
$ws.Range("J122").Value = 3324.75
$ws.Range("K122").Value = 10383.8568
$ws.Range("L122").Value = 9974.25
$ws.Range("M122").Value = -7933.856800000001
$ws.Range("N122").Value = -14874.25

# GSM!row 123 - Workplace Workout
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 8604.647000000001
$ws.Range("J123").Value = 8604.647000000001
$ws.Range("L123").Value = 8604.647000000001
$ws.Range("N123").Value = -13504.647

# GSM!row 132 - On Board for Lar
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3046.5881
$ws.Range("I132").Value = 2269.8572
$ws.Range("J132").Value = 6671.3335
$ws.Range("K132").Value = 6809.571599999999
$ws.Range("L132").Value = 20014.0005
$ws.Range("M132").Value = -4279.571599999999
$ws.Range("N132").Value = -25074.0005

# LTW!row 35 - No Risk, No Reward
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 2115.5
$ws.Range("I35").Value = 2115.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2115.5
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1779.5
$ws.Range("N35").ClearContents()

# LTW!row 40 - Best Served Toad
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3681
$ws.Range("I40").Value = 3600
$ws.Range("J40").Value = 3802.5
$ws.Range("K40").Value = 3600
$ws.Range("L40").Value = 3802.5
$ws.Range("M40").Value = -3464
$ws.Range("N40").Value = -4074.5

# LTW!row 122 - Hell on Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 50005556
$ws.Range("I122").Value = 83336670
$ws.Range("J122").Value = 33339998
$ws.Range("K122").Value = 250010010
$ws.Range("L122").Value = 100019994
$ws.Range("M122").Value = -250007560
$ws.Range("N122").Value = -100024894

# WVR!row 43 - Walk Softly and Carry a Big Halberd
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 45000
$ws.Range("I43").Value = 45000
$ws.Range("K43").Value = 45000
$ws.Range("M43").Value = -44851

# WVR!row 122 - Heavy Armoire
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 96155800
$ws.Range("I122").Value = 138890750
$ws.Range("J122").Value = 2164
$ws.Range("K122").Value = 416672250
$ws.Range("L122").Value = 6492
$ws.Range("M122").Value = -416669800
$ws.Range("N122").Value = -11392
